# Validating column data for empty/null values, DateOfBirth and IsActive
# columns on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IsActive (column C) was a "Y" shared-string placeholder for every row;
# replace it with real boolean TRUE/FALSE values so the column can be
# validated as a proper boolean column.
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $false

# Row 4 also gets an invalid DateOfBirth (plain text instead of a date) so
# the DateOfBirth column validation has a bad value to catch, alongside the
# existing row 5 which already has a blank Name to exercise empty/null
# checks.
$ws.Range("B4").Value = "fdfsdfsdf"
$ws.Range("C4").Value = $true

$ws.Range("C5").Value = $true
$ws.Range("C6").Value = $true

# Leave the selection on the helper column used while reviewing the data.
$ws.Range("F2:F6").Select() | Out-Null
